# Rebuild the sample rules sheet into a Drools-compliant decision table:
#   RuleSet / Import / Variables / Sequential header block, then a
#   RuleTable with NAME / CONDITION / CONDITION2 / CONDITION3 / ACTION / ACTION2
#   columns and three rule rows (YoungAdultRule / SeniorRule / PremiumRule).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start clean - remove all existing content/formatting from the old layout.
$ws.Cells.Clear()

function Set-TextValue {
    param($Range, [string]$Text)
    # Force genuinely-text storage even for values that look like numbers/
    # booleans (e.g. "18", "0.10", "TRUE") so Excel doesn't silently coerce
    # them into a numeric/boolean cell. A plain .Value assignment of such a
    # string gets auto-converted by Excel's type inference, and forcing a
    # text number-format leaves a stray (unused) style behind. Instead,
    # build a text-returning formula, then convert it to a static value via
    # copy / paste-values (xlPasteValues = -4163) - this keeps the cell on
    # the default/plain style with genuine shared-string storage.
    $escaped = $Text -replace '"', '""'
    $Range.Formula = '="" & "' + $escaped + '"'
    $Range.Copy()
    $Range.PasteSpecial(-4163)
}

# --- Header block (rows 1-4) ---
$ws.Range("A1").Value = "RuleSet"
$ws.Range("B1").Value = "CustomerRules"

$ws.Range("A2").Value = "Import"
$ws.Range("B2").Value = "com.example.model.Customer"

$ws.Range("A3").Value = "Variables"
$ws.Range("B3").Value = "Customer customer"

$ws.Range("A4").Value = "Sequential"
Set-TextValue $ws.Range("B4") "TRUE"

# Row 5 intentionally stays blank (matches the target layout's gap).

# --- RuleTable block (rows 6-11) ---
$ws.Range("A6").Value = "RuleTable"
$ws.Range("B6").Value = "CustomerDiscountRules"

$ws.Range("A7").Value = "NAME"
$ws.Range("B7").Value = "CONDITION"
$ws.Range("C7").Value = "CONDITION2"
$ws.Range("D7").Value = "CONDITION3"
$ws.Range("E7").Value = "ACTION"
$ws.Range("F7").Value = "ACTION2"

$ws.Range("B8").Value = "customer.getAge() >= `$param"
$ws.Range("C8").Value = 'customer.getStatus() == "$param"'
$ws.Range("D8").Value = "customer.isVip() == `$param"
$ws.Range("E8").Value = "customer.setDiscount(`$param);"
$ws.Range("F8").Value = 'customer.setStatus("$param");'

$ws.Range("A9").Value = "YoungAdultRule"
Set-TextValue $ws.Range("B9") "18"
$ws.Range("C9").Value = "ACTIVE"
Set-TextValue $ws.Range("E9") "0.05"

$ws.Range("A10").Value = "SeniorRule"
Set-TextValue $ws.Range("B10") "65"
$ws.Range("C10").Value = "SENIOR"
Set-TextValue $ws.Range("E10") "0.15"

$ws.Range("A11").Value = "PremiumRule"
Set-TextValue $ws.Range("B11") "25"
$ws.Range("C11").Value = "PREMIUM"
Set-TextValue $ws.Range("E11") "0.10"
